$d = $word.ActiveDocument

# Locate the field that holds the  m:'doc.html'.fromHTMLURI()  field code
# (today it is a classic  fldChar begin / instrText... / fldChar end  field)
# and rewrite it as plain literal-text runs wrapped in "{" / "}" braces, one
# run per original instrText chunk, exactly like the
# TokenIteratorFieldRewriterSplit parser expects the M2Doc template tokens
# to look.
$target = $null
foreach ($f in $d.Fields) {
    if ($f.Code.Text -match "fromHTMLURI") {
        $target = $f
    }
}

if ($target -eq $null) {
    throw "Could not find the doc.html fromHTMLURI field"
}

# Remember which paragraph owns the field so the replacement text can be
# inserted back into that very paragraph once the field is removed.
$paraIndex = 1
$i = 1
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Start -le $target.Code.Start -and $para.Range.End -ge $target.Code.End) {
        $paraIndex = $i
    }
    $i = $i + 1
}

# Delete the whole field: begin marker, instrText runs, end marker (the
# bookmark that used to sit in the middle of the field code is re-created
# below, in the same spot).
$target.Delete()

$para = $d.Paragraphs.Item($paraIndex)
$ins = $para.Range
$ins.Collapse(1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F"><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t>'</w:t></w:r><w:r><w:t>doc.html</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>'.fromHTMLURI()</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$null = $ins.InsertXML($xml)
